$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds numeric-looking values that must stay as TEXT
# (e.g. "94.954.46", "1.00", "0.999") -- mark the affected cells as Text
# format before assigning, so Excel does not reinterpret/round them as
# numbers. Only the cells whose Price actually changes are touched, so
# cells whose Price is untouched keep their original (General) style.
$ws.Range("D2:D6").NumberFormat = "@"
$ws.Range("D8:D11").NumberFormat = "@"
$ws.Range("D13:D16").NumberFormat = "@"
$ws.Range("D18:D33").NumberFormat = "@"
$ws.Range("D35:D40").NumberFormat = "@"
$ws.Range("D43:D51").NumberFormat = "@"

$ws.Range("D2").Value = "94.954.46"
$ws.Range("E2").Value = "  -1.44%  "

$ws.Range("D3").Value = "3.461.85"
$ws.Range("E3").Value = "  +4.18%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "239.79"
$ws.Range("E5").Value = "  -3.47%  "

$ws.Range("D6").Value = "644.47"
$ws.Range("E6").Value = "  -0.96%  "

$ws.Range("E7").Value = "  +6.38%  "

$ws.Range("D8").Value = "0.405"
$ws.Range("E8").Value = "  -3.12%  "

$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").Value = "1.01"
$ws.Range("E10").Value = "  +2.80%  "

$ws.Range("D11").Value = "3.462.13"
$ws.Range("E11").Value = "  +4.24%  "

$ws.Range("E12").Value = "  -3.40%  "

$ws.Range("D13").Value = "41.89"
$ws.Range("E13").Value = "  +4.68%  "

$ws.Range("D14").Value = "6.14"
$ws.Range("E14").Value = "  +1.59%  "

$ws.Range("D15").Value = "94.694.67"
$ws.Range("E15").Value = "  -1.42%  "

$ws.Range("D16").Value = "4.105.78"
$ws.Range("E16").Value = "  +4.33%  "

$ws.Range("E17").Value = "  +2.83%  "

$ws.Range("D18").Value = "8.52"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").Value = "3.453.21"
$ws.Range("E19").Value = "  +3.69%  "

$ws.Range("D20").Value = "17.92"
$ws.Range("E20").Value = "  +5.57%  "

$ws.Range("D21").Value = "11.45"
$ws.Range("E21").Value = "  +9.71%  "

$ws.Range("D22").Value = "0.511"
$ws.Range("E22").Value = "  -4.25%  "

$ws.Range("D23").Value = "503.61"
$ws.Range("E23").Value = "  +0.25%  "

$ws.Range("D24").Value = "3.19"
$ws.Range("E24").Value = "  -4.82%  "

$ws.Range("D25").Value = "0.0000193"
$ws.Range("E25").Value = "  -1.62%  "

$ws.Range("D26").Value = "6.51"
$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("D27").Value = "91.90"
$ws.Range("E27").Value = "  -3.97%  "

$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "12.16"
$ws.Range("E28").Value = "  +1.50%  "

$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.637.52"
$ws.Range("E29").Value = "  +4.02%  "

$ws.Range("D30").Value = "11.76"
$ws.Range("E30").Value = "  +7.66%  "

$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("D32").Value = "2.75"
$ws.Range("E32").Value = "  +11.87%  "

$ws.Range("D33").Value = "0.138"
$ws.Range("E33").Value = "  -3.13%  "

$ws.Range("E34").Value = "  -1.66%  "

$ws.Range("D35").Value = "31.02"
$ws.Range("E35").Value = "  +11.37%  "

$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").Value = "0.567"
$ws.Range("E37").Value = "  +4.50%  "

$ws.Range("D38").Value = "7.74"
$ws.Range("E38").Value = "  +2.21%  "

$ws.Range("D39").Value = "1.45"
$ws.Range("E39").Value = "  -0.55%  "

$ws.Range("D40").Value = "525.35"
$ws.Range("E40").Value = "  +4.44%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("E42").Value = "  +0.22%  "

$ws.Range("D43").Value = "0.925"
$ws.Range("E43").Value = "  +12.11%  "

$ws.Range("D44").Value = "24.08"
$ws.Range("E44").Value = "  -1.09%  "

$ws.Range("D45").Value = "5.69"
$ws.Range("E45").Value = "  +4.22%  "

$ws.Range("D46").Value = "1.71"
$ws.Range("E46").Value = "  +3.22%  "

$ws.Range("D47").Value = "0.0417"
$ws.Range("E47").Value = "  -2.06%  "

$ws.Range("D48").Value = "3.50"
$ws.Range("E48").Value = "  -3.83%  "

$ws.Range("D49").Value = "2.16"
$ws.Range("E49").Value = "  +10.25%  "

$ws.Range("D50").Value = "53.38"
$ws.Range("E50").Value = "  +0.59%  "

$ws.Range("D51").Value = "3.21"
$ws.Range("E51").Value = "  +3.34%  "
